# Rename the four worksheets to their new, more descriptive names.
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("msq").Name     = "msq_node"
$wb.Worksheets.Item("kfk").Name     = "kfk_node"
$wb.Worksheets.Item("msq_msq").Name = "msq_msq_rule"
$wb.Worksheets.Item("msq_kfk").Name = "msq_kfk_rule"

# Restore/move each sheet's cursor position. Selecting a range both sets the
# sheet's activeCell/selection and activates that sheet (and tab), so the
# final Select() below (on msq_kfk_rule) is what ends up as the workbook's
# active tab - matching the target state.
$wb.Worksheets.Item("msq_node").Range("B43").Select()     | Out-Null
$wb.Worksheets.Item("kfk_node").Range("C94").Select()      | Out-Null
$wb.Worksheets.Item("msq_msq_rule").Range("B123").Select() | Out-Null
$wb.Worksheets.Item("msq_kfk_rule").Range("E265").Select() | Out-Null
